# Applies the "spell/gram proofing markup + expanded instructions" edit
# described by the diff. We rebuild each touched paragraph's contents via
# Range.InsertXML (which replaces that paragraph's contents in place,
# letting us place w:proofErr markers exactly where Word would) and we
# insert brand-new paragraphs with Range.InsertParagraphAfter() + InsertXML.
#
# We walk the document from the LAST touched paragraph back to the FIRST,
# so that paragraph insertions (which shift every later index) never
# invalidate an index we still need to use.

$d = $word.ActiveDocument
$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------
# Paragraph 18: "5.-"  -> move the _GoBack bookmark here, after the run.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(18)
$xml = "<w:p $W>" +
       "<w:r><w:t>5.-</w:t></w:r>" +
       "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# Paragraph 15: "4.-After the first CMake configuration for nscale is
# successful check NS_SEGMENT  as ON and reconfigure the project" gains
# proofing marks plus a brand-new trailing sentence.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(15)
$xml = "<w:p $W>" +
       "<w:r><w:t xml:space='preserve'>4.-After the first </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>CMake</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> configuration for </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>nscale</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> is successful</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> check NS_</w:t></w:r>" +
       "<w:proofErr w:type='gramStart'/><w:r><w:t>SEGMENT  as</w:t></w:r><w:proofErr w:type='gramEnd'/>" +
       "<w:r><w:t xml:space='preserve'> ON and reconfigure the project</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'>. New options will show up. If you want to use Open CV and </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>nscale</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> CUDA-GPU functionality you need to make sure to check the USE_CUDA option.</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# Paragraph 13 used to read "ZLIB: use <openCV fodler>\3rdparty\zlib".
# It becomes the new "If nscale Cmake configuration cannot find Zlib..."
# paragraph, and the real (now proofed) ZLIB line is re-created as two
# brand-new paragraphs right after it (a blank spacer + the ZLIB text).
# ---------------------------------------------------------------------
$p = $d.Paragraphs(13)
$xml = "<w:p $W>" +
       "<w:r><w:t>If</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>nscale</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>Cmake</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> configuration cannot find </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>Zlib</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> then</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> you need to indicate the following directories:</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# One new blank spacer paragraph, then one new paragraph for the ZLIB line.
$p.Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs(14).Range.InsertXML("<w:p $W/>")
$d.Paragraphs(14).Range.InsertParagraphAfter() | Out-Null
$xml = "<w:p $W>" +
       "<w:r><w:t>ZLIB: use &lt;</w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>openCV</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>fodler</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t>&gt;\3rdparty\</w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>zlib</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "</w:p>"
$d.Paragraphs(15).Range.InsertXML($xml)

# ---------------------------------------------------------------------
# Paragraph 11: "3.-In nscale Cmake you need to indicate the following
# directories:" becomes "3.-Open CMake-gui for nscale project."
# ---------------------------------------------------------------------
$p = $d.Paragraphs(11)
$xml = "<w:p $W>" +
       "<w:proofErr w:type='gramStart'/><w:r><w:t>3.-</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'>Open </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>CMake-gui</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> for </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>nscale</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> project.</w:t></w:r>" +
       "<w:proofErr w:type='gramEnd'/>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# Paragraph 9: "PS: for Visual Studio 2012 add ;%OPENCV_DIR%\x64\vc11\bin"
# loses the _GoBack bookmark (moved to paragraph 18 above) and gains
# proofing marks around "add ;%".
# ---------------------------------------------------------------------
$p = $d.Paragraphs(9)
$xml = "<w:p $W>" +
       "<w:r><w:t>PS: for Visual Studio 2012</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
       "<w:proofErr w:type='gramStart'/><w:r><w:t>add ;%</w:t></w:r><w:proofErr w:type='gramEnd'/>" +
       "<w:r><w:t>OPENCV_DIR%\x64</w:t></w:r>" +
       "<w:r><w:t>\vc11\bin</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# Paragraph 8: "PS: for Visual Studio 2013 add ;%OPENCV_DIR%\x64\vc12\bin"
# gains proofing marks around "add ;%" (trailing x64 / \vc12\bin runs
# are preserved as-is).
# ---------------------------------------------------------------------
$p = $d.Paragraphs(8)
$xml = "<w:p $W>" +
       "<w:r><w:t xml:space='preserve'>PS: for Visual Studio 2013 </w:t></w:r>" +
       "<w:proofErr w:type='gramStart'/><w:r><w:t>add ;%</w:t></w:r><w:proofErr w:type='gramEnd'/>" +
       "<w:r><w:t>OPENCV_DIR%\</w:t></w:r>" +
       "<w:r><w:t>x64</w:t></w:r>" +
       "<w:r><w:t>\vc12\bin</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# Paragraph 5: "OPENCV_DIR = <openCV fodler>\build" gains proofing marks.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(5)
$xml = "<w:p $W>" +
       "<w:r><w:t>OPENCV_DIR = &lt;</w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>openCV</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>fodler</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t>&gt;\build</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# Paragraph 2: "1.-Download and install Opencv 2.4.10 from here ..."
# gains proofing marks around "Opencv".
# ---------------------------------------------------------------------
$p = $d.Paragraphs(2)
$xml = "<w:p $W>" +
       "<w:r><w:t xml:space='preserve'>1.-Download and install </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>Opencv</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> 2.4.10 from here http://sourceforge.net/projects/opencvlibrary/files/latest/download?source=files</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# Paragraph 1: "INSTALL AND COMPILE nscale library in Windows" gains
# proofing marks around "nscale".
# ---------------------------------------------------------------------
$p = $d.Paragraphs(1)
$xml = "<w:p $W>" +
       "<w:r><w:t xml:space='preserve'>INSTALL AND COMPILE </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/><w:r><w:t>nscale</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> library in Windows</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

Write-Output "Done. Final paragraph count: $($d.Paragraphs.Count)"
